$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1647
$ws.Cells.Item(4, 6).Value = 9458
$ws.Cells.Item(5, 6).Value = 734
$ws.Cells.Item(6, 6).Value = 618
$ws.Cells.Item(7, 6).Value = 212
$ws.Cells.Item(8, 6).Value = 323
$ws.Cells.Item(11, 6).Value = 1626
$ws.Cells.Item(12, 6).Value = 1403
$ws.Cells.Item(15, 6).Value = 1467
$ws.Cells.Item(16, 6).Value = 123
$ws.Cells.Item(17, 6).Value = 300
$ws.Cells.Item(19, 6).Value = 137
$ws.Cells.Item(20, 6).Value = 82
$ws.Cells.Item(21, 6).Value = 376
$ws.Cells.Item(23, 6).Value = 93
$ws.Cells.Item(28, 6).Value = 254
$ws.Cells.Item(33, 6).Value = 165
$ws.Cells.Item(35, 6).Value = 179
$ws.Cells.Item(38, 6).Value = 234
$ws.Cells.Item(39, 6).Value = 604
$ws.Cells.Item(42, 6).Value = 734
$ws.Cells.Item(45, 6).Value = 314

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 67
$ws.Cells.Item(15, 6).Value = 46
$ws.Cells.Item(18, 6).Value = 966
$ws.Cells.Item(19, 6).Value = 31
$ws.Cells.Item(20, 6).Value = 1061
$ws.Cells.Item(21, 6).Value = 280
$ws.Cells.Item(22, 6).Value = 654
$ws.Cells.Item(23, 6).Value = 11
$ws.Cells.Item(25, 6).Value = 308
$ws.Cells.Item(31, 6).Value = 169
$ws.Cells.Item(35, 6).Value = 113
$ws.Cells.Item(37, 6).Value = 16
$ws.Cells.Item(38, 6).Value = 24

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 771
$ws.Cells.Item(7, 6).Value = 2358
$ws.Cells.Item(8, 6).Value = 3572
$ws.Cells.Item(9, 6).Value = 7
$ws.Cells.Item(11, 6).Value = 63
$ws.Cells.Item(12, 6).Value = 102

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1647
$ws.Cells.Item(3, 6).Value = 9458
$ws.Cells.Item(6, 6).Value = 3572
$ws.Cells.Item(7, 6).Value = 734
$ws.Cells.Item(8, 6).Value = 63
$ws.Cells.Item(9, 6).Value = 63
$ws.Cells.Item(10, 6).Value = 618
$ws.Cells.Item(11, 6).Value = 212
$ws.Cells.Item(12, 6).Value = 323
$ws.Cells.Item(14, 6).Value = 1403
$ws.Cells.Item(16, 6).Value = 102
$ws.Cells.Item(17, 6).Value = 102
$ws.Cells.Item(18, 6).Value = 1467
$ws.Cells.Item(19, 6).Value = 300
$ws.Cells.Item(21, 6).Value = 137
$ws.Cells.Item(23, 6).Value = 93
$ws.Cells.Item(26, 6).Value = 46
$ws.Cells.Item(29, 6).Value = 31
$ws.Cells.Item(31, 6).Value = 254
$ws.Cells.Item(32, 6).Value = 1061
$ws.Cells.Item(33, 6).Value = 280
$ws.Cells.Item(36, 6).Value = 11
$ws.Cells.Item(38, 6).Value = 165
$ws.Cells.Item(39, 6).Value = 308
$ws.Cells.Item(40, 6).Value = 308
$ws.Cells.Item(43, 6).Value = 234
$ws.Cells.Item(45, 6).Value = 604
$ws.Cells.Item(47, 6).Value = 734
$ws.Cells.Item(49, 6).Value = 113
$ws.Cells.Item(50, 6).Value = 314
$ws.Cells.Item(52, 6).Value = 16
